# Basic setup for spread model setup
# Populate rows 2-5 with updated computed values and team names, and add a new row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 234.5
$ws.Range("E2").Value = 17253288.5
$ws.Range("F2").Value = 115.3205128205128
$ws.Range("G2").Value = 13.5
$ws.Range("H2").Value = "SanAntonio"
$ws.Range("I2").Value = "Boston"
$ws.Range("J2").Value = 0.5462304409672831
$ws.Range("K2").Value = 99.5051282051282
$ws.Range("L2").Value = 115.3269230769231
$ws.Range("M2").Value = 116.55
$ws.Range("N2").Value = 76.61538461538463
$ws.Range("O2").Value = 0.4154102564102563
$ws.Range("P2").Value = 0.5860769230769232
$ws.Range("Q2").Value = 0.2534230769230769
$ws.Range("R2").Value = 12.70128205128205
$ws.Range("S2").Value = 11.8525641025641
$ws.Range("T2").Value = 0.1952115384615385
$ws.Range("U2").Value = 1.013361272587986
$ws.Range("V2").Value = 0.9974903451022887
$ws.Range("W2").Value = 11.12508988895773

# --- Row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 238.5
$ws.Range("E3").Value = 19299819
$ws.Range("F3").Value = 116.1529080675422
$ws.Range("G3").Value = 1.5
$ws.Range("H3").Value = "Chicago"
$ws.Range("I3").Value = "Utah"
$ws.Range("J3").Value = 0.5185897435897436
$ws.Range("K3").Value = 99.31041275797372
$ws.Range("L3").Value = 115.8054721701063
$ws.Range("M3").Value = 115.7099749843652
$ws.Range("N3").Value = 76.27911194496559
$ws.Range("O3").Value = 0.3890678549093183
$ws.Range("P3").Value = 0.5898373983739836
$ws.Range("Q3").Value = 0.2630712945590995
$ws.Range("R3").Value = 12.4750469043152
$ws.Range("S3").Value = 12.53164477798624
$ws.Range("T3").Value = 0.2111558786741714
$ws.Range("U3").Value = 1.020675817816715
$ws.Range("V3").Value = 1.076525518682472
$ws.Range("W3").Value = 10.68952267008545

# --- Row 4 ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 225.5
$ws.Range("E4").Value = 8276298.5
$ws.Range("F4").Value = 114.6025641025641
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = "Dallas"
$ws.Range("I4").Value = "NewOrleans"
$ws.Range("J4").Value = 0.5641025641025641
$ws.Range("K4").Value = 97.23846153846154
$ws.Range("L4").Value = 116.2141025641026
$ws.Range("M4").Value = 113.2615384615385
$ws.Range("N4").Value = 77.02051282051281
$ws.Range("O4").Value = 0.4216666666666667
$ws.Range("P4").Value = 0.5883589743589743
$ws.Range("Q4").Value = 0.3035897435897436
$ws.Range("R4").Value = 11.9051282051282
$ws.Range("S4").Value = 12.92564102564102
$ws.Range("T4").Value = 0.2227948717948718
$ws.Range("U4").Value = 1.007052408634131
$ws.Range("V4").Value = 0.9661064756566271
$ws.Range("W4").Value = 10.77742982739587

# --- Row 5 ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 231.5
$ws.Range("E5").Value = 27155624
$ws.Range("F5").Value = 113.7307692307692
$ws.Range("G5").Value = 6.5
$ws.Range("H5").Value = "GoldenState"
$ws.Range("I5").Value = "Orlando"
$ws.Range("J5").Value = 0.5597439544807966
$ws.Range("K5").Value = 99.80512820512823
$ws.Range("L5").Value = 113.2346153846154
$ws.Range("M5").Value = 115.4423076923077
$ws.Range("N5").Value = 76.49871794871795
$ws.Range("O5").Value = 0.4217179487179488
$ws.Range("P5").Value = 0.5855384615384613
$ws.Range("Q5").Value = 0.2613974358974359
$ws.Range("R5").Value = 13.76794871794872
$ws.Range("S5").Value = 11.97435897435897
$ws.Range("T5").Value = 0.2149807692307692
$ws.Range("U5").Value = 0.9993916452615925
$ws.Range("V5").Value = 1.055825170109766
$ws.Range("W5").Value = 11.02751638816982

# --- Row 6 (new row) ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 241.5
$ws.Range("E6").Value = 31473048
$ws.Range("F6").Value = 117.1389466389466
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = "Sacramento"
$ws.Range("I6").Value = "LALakers"
$ws.Range("J6").Value = 0.517094017094017
$ws.Range("K6").Value = 101.4100831600831
$ws.Range("L6").Value = 114.9733194733195
$ws.Range("M6").Value = 114.9518018018018
$ws.Range("N6").Value = 77.20772695772695
$ws.Range("O6").Value = 0.3840717255717255
$ws.Range("P6").Value = 0.5918808038808039
$ws.Range("Q6").Value = 0.2856933471933473
$ws.Range("R6").Value = 12.595841995842
$ws.Range("S6").Value = 11.55696465696466
$ws.Range("T6").Value = 0.2125571725571726
$ws.Range("U6").Value = 1.029340480131341
$ws.Range("V6").Value = 1.0033256993911
$ws.Range("W6").Value = 11.30656842760246

# Match formatting of the existing "index" column cells (A2:A5) for the new A6 cell
# (bordered/centered/bold style), by copying format from A5 down to A6.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
